$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary cell updates ---
$ws.Range("E11").Value = 753999
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 12

# --- Copy "closing" border formatting from the last existing data row (36) to what will become the new last data row (30) ---
$ws.Range("B36:J36").Copy()
$ws.Range("B30:J30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows: updated employee account-statement entries ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1007857666"
$ws.Range("D16").Value = "LEONARDO JOSE MENDEZ AISLANT"
$ws.Range("E16").Value = "2409"
$ws.Range("F16").Value = 34666
$ws.Range("G16").Value = 1300000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1007857666"
$ws.Range("D17").Value = "LEONARDO JOSE MENDEZ AISLANT"
$ws.Range("E17").Value = "2410"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1300000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "92131117"
$ws.Range("D18").Value = "SAMID ANTONIO RUIZ MERIÝO"
$ws.Range("E18").Value = "2411"
$ws.Range("F18").Value = 43333
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1007857666"
$ws.Range("D19").Value = "LEONARDO JOSE MENDEZ AISLANT"
$ws.Range("E19").Value = "2411"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1065875439"
$ws.Range("D20").Value = "JAN ESLEIDER RINALDY QUINTERO"
$ws.Range("E20").Value = "2412"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1423500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "92131117"
$ws.Range("D21").Value = "SAMID ANTONIO RUIZ MERIÝO"
$ws.Range("E21").Value = "2412"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 1423500

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1007857666"
$ws.Range("D22").Value = "LEONARDO JOSE MENDEZ AISLANT"
$ws.Range("E22").Value = "2412"
$ws.Range("F22").Value = 52000
$ws.Range("G22").Value = 1300000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1007857666"
$ws.Range("D23").Value = "LEONARDO JOSE MENDEZ AISLANT"
$ws.Range("E23").Value = "2501"
$ws.Range("F23").Value = 52000
$ws.Range("G23").Value = 1300000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1007857666"
$ws.Range("D24").Value = "LEONARDO JOSE MENDEZ AISLANT"
$ws.Range("E24").Value = "2502"
$ws.Range("F24").Value = 52000
$ws.Range("G24").Value = 1300000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1007857666"
$ws.Range("D25").Value = "LEONARDO JOSE MENDEZ AISLANT"
$ws.Range("E25").Value = "2503"
$ws.Range("F25").Value = 52000
$ws.Range("G25").Value = 1300000

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1007857666"
$ws.Range("D26").Value = "LEONARDO JOSE MENDEZ AISLANT"
$ws.Range("E26").Value = "2504"
$ws.Range("F26").Value = 52000
$ws.Range("G26").Value = 1300000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1007857666"
$ws.Range("D27").Value = "LEONARDO JOSE MENDEZ AISLANT"
$ws.Range("E27").Value = "2505"
$ws.Range("F27").Value = 52000
$ws.Range("G27").Value = 1300000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1007857666"
$ws.Range("D28").Value = "LEONARDO JOSE MENDEZ AISLANT"
$ws.Range("E28").Value = "2506"
$ws.Range("F28").Value = 52000
$ws.Range("G28").Value = 1300000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1007857666"
$ws.Range("D29").Value = "LEONARDO JOSE MENDEZ AISLANT"
$ws.Range("E29").Value = "2507"
$ws.Range("F29").Value = 52000
$ws.Range("G29").Value = 1300000

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "1007857666"
$ws.Range("D30").Value = "LEONARDO JOSE MENDEZ AISLANT"
$ws.Range("E30").Value = "2508"
$ws.Range("F30").Value = 52000
$ws.Range("G30").Value = 1300000

# --- Remove the now-obsolete trailing data rows (31:36); this shifts the signature-block rows (formerly 41/42) up to 35/36 ---
$ws.Rows("31:36").Delete()

# --- Column D (Nombre Trabajador) best-fit width shrank now that the longest name is shorter ---
$ws.Columns("D").ColumnWidth = 31.53
